$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1500
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -1850

$ws.Range("H74").Value = 6250
$ws.Range("I74").Value = 6250
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 6250
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -5314

$ws.Range("H77").Value = 6250
$ws.Range("I77").Value = 6250
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 31250
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -26570

$ws.Range("H100").Value = 1724.4
$ws.Range("I100").Value = 1487.9333
$ws.Range("J100").Value = 2433.8
$ws.Range("K100").Value = 1487.9333
$ws.Range("L100").Value = 2433.8
$ws.Range("M100").Value = -946.9332999999999
$ws.Range("N100").Value = -3515.8

$ws.Range("H112").Value = 3133.4285
$ws.Range("I112").Value = 2872.5
$ws.Range("J112").Value = 3329.125
$ws.Range("K112").Value = 8617.5
$ws.Range("L112").Value = 9987.375
$ws.Range("M112").Value = -7509.5
$ws.Range("N112").Value = -12203.375

$ws.Range("H131").Value = 1446.8462
$ws.Range("I131").Value = 1442.8334
$ws.Range("J131").Value = 1495
$ws.Range("K131").Value = 4328.5002
$ws.Range("L131").Value = 4485
$ws.Range("M131").Value = 711.4997999999996
$ws.Range("N131").Value = -14565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2990.484
$ws.Range("I61").Value = 2790.1667
$ws.Range("J61").Value = 9000
$ws.Range("K61").Value = 2790.1667
$ws.Range("L61").Value = 9000
$ws.Range("M61").Value = -2578.1667
$ws.Range("N61").Value = -9424

$ws.Range("H88").Value = 2409.0557
$ws.Range("I88").Value = 2035.6923
$ws.Range("J88").Value = 3379.8
$ws.Range("K88").Value = 2035.6923
$ws.Range("L88").Value = 3379.8
$ws.Range("M88").Value = -1629.6923
$ws.Range("N88").Value = -4191.8

$ws.Range("H91").Value = 2409.0557
$ws.Range("I91").Value = 2035.6923
$ws.Range("J91").Value = 3379.8
$ws.Range("K91").Value = 2035.6923
$ws.Range("L91").Value = 3379.8
$ws.Range("M91").Value = -631.6922999999999
$ws.Range("N91").Value = -6187.8

$ws.Range("H132").Value = 3818.3635
$ws.Range("I132").Value = 2318.2969
$ws.Range("J132").Value = 7818.5415
$ws.Range("K132").Value = 6954.8907
$ws.Range("L132").Value = 23455.6245
$ws.Range("M132").Value = -4424.8907
$ws.Range("N132").Value = -28515.6245

$ws.Range("H136").Value = 2990.484
$ws.Range("I136").Value = 2790.1667
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 8370.500100000001
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -5820.500100000001
$ws.Range("N136").Value = -32100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 70660
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 70660
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 70660
$ws.Range("N13").Value = -70996

$ws.Range("H86").Value = 1905.4286
$ws.Range("I86").Value = 1973
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 1973
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -850
$ws.Range("N86").Value = -3746

$ws.Range("H89").Value = 1905.4286
$ws.Range("I89").Value = 1973
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 9865
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -4249
$ws.Range("N89").Value = -18732

$ws.Range("H94").Value = 7950.8
$ws.Range("I94").Value = 7506.6
$ws.Range("J94").Value = 9283.4
$ws.Range("K94").Value = 7506.6
$ws.Range("L94").Value = 9283.4
$ws.Range("M94").Value = -7055.6
$ws.Range("N94").Value = -10185.4

$ws.Range("H100").Value = 23500
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 23500
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 23500
$ws.Range("N100").Value = -25664
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5063
$ws.Range("I58").Value = 3348.7896
$ws.Range("J58").Value = 9134.25
$ws.Range("K58").Value = 3348.7896
$ws.Range("L58").Value = 9134.25
$ws.Range("M58").Value = -3145.7896
$ws.Range("N58").Value = -9540.25

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H80").Value = 43128
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 43128
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 43128
$ws.Range("N80").Value = -45374

$ws.Range("H83").Value = 43128
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 43128
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 129384
$ws.Range("N83").Value = -140616

$ws.Range("H99").Value = 12057.429
$ws.Range("I99").Value = 18974.916
$ws.Range("J99").Value = 2834.111
$ws.Range("K99").Value = 18974.916
$ws.Range("L99").Value = 2834.111
$ws.Range("M99").Value = -17476.916
$ws.Range("N99").Value = -5830.111

$ws.Range("H126").Value = 12057.429
$ws.Range("I126").Value = 18974.916
$ws.Range("J126").Value = 2834.111
$ws.Range("K126").Value = 56924.74800000001
$ws.Range("L126").Value = 8502.332999999999
$ws.Range("M126").Value = -54454.74800000001
$ws.Range("N126").Value = -13442.333

$ws.Range("H132").Value = 4203.6665
$ws.Range("I132").Value = 4188.4287
$ws.Range("J132").Value = 4257
$ws.Range("K132").Value = 12565.2861
$ws.Range("L132").Value = 12771
$ws.Range("M132").Value = -10035.2861
$ws.Range("N132").Value = -17831

$ws.Range("H134").Value = 3023.077
$ws.Range("I134").Value = 3066.75
$ws.Range("J134").Value = 2499
$ws.Range("K134").Value = 9200.25
$ws.Range("L134").Value = 7497
$ws.Range("M134").Value = -6665.25
$ws.Range("N134").Value = -12567

$ws.Range("H136").Value = 5063
$ws.Range("I136").Value = 3348.7896
$ws.Range("J136").Value = 9134.25
$ws.Range("K136").Value = 10046.3688
$ws.Range("L136").Value = 27402.75
$ws.Range("M136").Value = -7496.3688
$ws.Range("N136").Value = -32502.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2421.1333
$ws.Range("I34").Value = 1603.1666
$ws.Range("J34").Value = 2966.4443
$ws.Range("K34").Value = 4809.4998
$ws.Range("L34").Value = 8899.332900000001
$ws.Range("M34").Value = -4725.4998
$ws.Range("N34").Value = -9067.332900000001

$ws.Range("H110").Value = 16372
$ws.Range("I110").Value = 10325.333
$ws.Range("J110").Value = 20000
$ws.Range("K110").Value = 30975.999
$ws.Range("L110").Value = 60000
$ws.Range("M110").Value = -26885.999
$ws.Range("N110").Value = -68180

$ws.Range("H119").Value = 9621.25
$ws.Range("I119").Value = 2207.8572
$ws.Range("J119").Value = 20000
$ws.Range("K119").Value = 6623.571599999999
$ws.Range("L119").Value = 60000
$ws.Range("M119").Value = -1785.571599999999
$ws.Range("N119").Value = -69676

$ws.Range("H121").Value = 15742.5
$ws.Range("I121").Value = 259.75
$ws.Range("J121").Value = 28128.7
$ws.Range("K121").Value = 779.25
$ws.Range("L121").Value = 84386.10000000001
$ws.Range("M121").Value = 530.75
$ws.Range("N121").Value = -87006.10000000001

$ws.Range("H131").Value = 2264.1755
$ws.Range("I131").Value = 823.875
$ws.Range("J131").Value = 2499.3264
$ws.Range("K131").Value = 2471.625
$ws.Range("L131").Value = 7497.9792
$ws.Range("M131").Value = 2568.375
$ws.Range("N131").Value = -17577.9792

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 46933.332
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 46933.332
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 46933.332
$ws.Range("N45").Value = -48051.332

$ws.Range("H80").Value = 3989.5454
$ws.Range("I80").Value = 3465
$ws.Range("J80").Value = 6350
$ws.Range("K80").Value = 3465
$ws.Range("L80").Value = 6350
$ws.Range("M80").Value = -2467
$ws.Range("N80").Value = -8346

$ws.Range("H83").Value = 3989.5454
$ws.Range("I83").Value = 3465
$ws.Range("J83").Value = 6350
$ws.Range("K83").Value = 17325
$ws.Range("L83").Value = 31750
$ws.Range("M83").Value = -12333
$ws.Range("N83").Value = -41734

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3393.7
$ws.Range("I22").Value = 2193.5
$ws.Range("J22").Value = 3527.0557
$ws.Range("K22").Value = 2193.5
$ws.Range("L22").Value = 3527.0557
$ws.Range("M22").Value = -1898.5
$ws.Range("N22").Value = -4117.0557

$ws.Range("H27").Value = 3393.7
$ws.Range("I27").Value = 2193.5
$ws.Range("J27").Value = 3527.0557
$ws.Range("K27").Value = 2193.5
$ws.Range("L27").Value = 3527.0557
$ws.Range("M27").Value = -2086.5
$ws.Range("N27").Value = -3741.0557

$ws.Range("H46").Value = 4650.3477
$ws.Range("I46").Value = 2661.5
$ws.Range("J46").Value = 5352.294
$ws.Range("K46").Value = 2661.5
$ws.Range("L46").Value = 5352.294
$ws.Range("M46").Value = -2473.5
$ws.Range("N46").Value = -5728.294

$ws.Range("H61").Value = 10168.794
$ws.Range("I61").Value = 9435
$ws.Range("J61").Value = 12553.625
$ws.Range("K61").Value = 9435
$ws.Range("L61").Value = 12553.625
$ws.Range("M61").Value = -9233
$ws.Range("N61").Value = -12957.625

$ws.Range("H99").Value = 41722.5
$ws.Range("I99").Value = 43296.668
$ws.Range("J99").Value = 37000
$ws.Range("K99").Value = 43296.668
$ws.Range("L99").Value = 37000
$ws.Range("M99").Value = -40301.668
$ws.Range("N99").Value = -42990

$ws.Range("H113").Value = 10168.794
$ws.Range("I113").Value = 9435
$ws.Range("J113").Value = 12553.625
$ws.Range("K113").Value = 9435
$ws.Range("L113").Value = 12553.625
$ws.Range("M113").Value = -7265
$ws.Range("N113").Value = -16893.625

$ws.Range("H122").Value = 6025.5
$ws.Range("I122").Value = 4042
$ws.Range("J122").Value = 9992.5
$ws.Range("K122").Value = 12126
$ws.Range("L122").Value = 29977.5
$ws.Range("M122").Value = -9676
$ws.Range("N122").Value = -34877.5

$ws.Range("H132").Value = 2907.2432
$ws.Range("I132").Value = 2444.2415
$ws.Range("J132").Value = 4585.625
$ws.Range("K132").Value = 7332.7245
$ws.Range("L132").Value = 13756.875
$ws.Range("M132").Value = -4802.7245
$ws.Range("N132").Value = -18816.875

$ws.Range("H136").Value = 3196.8948
$ws.Range("I136").Value = 2865.875
$ws.Range("J136").Value = 4962.3335
$ws.Range("K136").Value = 8597.625
$ws.Range("L136").Value = 14887.0005
$ws.Range("M136").Value = -6047.625
$ws.Range("N136").Value = -19987.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 27972.666
$ws.Range("I41").Value = 7499
$ws.Range("J41").Value = 29833.908
$ws.Range("K41").Value = 7499
$ws.Range("L41").Value = 29833.908
$ws.Range("M41").Value = -7109
$ws.Range("N41").Value = -30613.908

$ws.Range("H62").Value = 8713.482
$ws.Range("I62").Value = 12648.625
$ws.Range("J62").Value = 7214.381
$ws.Range("K62").Value = 12648.625
$ws.Range("L62").Value = 7214.381
$ws.Range("M62").Value = -12024.625
$ws.Range("N62").Value = -8462.381000000001

$ws.Range("H65").Value = 8713.482
$ws.Range("I65").Value = 12648.625
$ws.Range("J65").Value = 7214.381
$ws.Range("K65").Value = 63243.125
$ws.Range("L65").Value = 36071.905
$ws.Range("M65").Value = -60123.125
$ws.Range("N65").Value = -42311.905

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H136").Value = 1721.2
$ws.Range("I136").Value = 1469.5769
$ws.Range("J136").Value = 2448.111
$ws.Range("K136").Value = 4408.7307
$ws.Range("L136").Value = 7344.333
$ws.Range("M136").Value = -1858.7307
$ws.Range("N136").Value = -12444.333
